# implemented new partition method - link process method
#
# Sheet1 holds an eGRID-style emission factor table (solid/liq/gas mg/mwh
# per region). The new partition method recomputes the per-region numbers;
# rows that have no contribution under the new method are left blank
# (only the region-name cell in column A remains) instead of storing
# explicit zeros, and the header row's previously-zeroed "total" row
# (row 2) is cleared out entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 had no region label and only zeros - clear it completely so the
# row drops out of the sheet's used data entirely.
$ws.Range("B2:D2").ClearContents()

# Regions with no measured contribution under the new method: keep the
# region-name cell (column A) but blank out the B:D metrics.
$ws.Range("B3:D3").ClearContents()
$ws.Range("B4:D4").ClearContents()
$ws.Range("B6:D6").ClearContents()
$ws.Range("B9:D9").ClearContents()
$ws.Range("B10:D10").ClearContents()
$ws.Range("B13:D13").ClearContents()
$ws.Range("B15:D15").ClearContents()
$ws.Range("B16:D16").ClearContents()

# Recomputed values for the remaining regions (new partition / link
# process method).
$ws.Range("B5").Value = 1121.6475366560899
$ws.Range("C5").Value = 0.45866400042399974
$ws.Range("D5").Value = 12.735575250268212

$ws.Range("B7").Value = 2397.4469980336244
$ws.Range("C7").Value = 0.16629701509286393
$ws.Range("D7").Value = 3.4675820974288225

$ws.Range("B8").Value = 2804.8673341435751
$ws.Range("C8").Value = 0.24632350762520833
$ws.Range("D8").Value = 5.2791067366499052

$ws.Range("B11").Value = 1046.7037338972473
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 122.69366375369066

$ws.Range("B12").Value = 1271.3509444178469
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 110.29082062794942

$ws.Range("B14").Value = 1105.688290913879
$ws.Range("C14").Value = 0.014137397268505738
$ws.Range("D14").Value = 33.298787647873297

$ws.Range("B17").Value = 5618.8865103340195
$ws.Range("C17").Value = 0.3987294610427321
$ws.Range("D17").Value = 7.3607681001160614

$ws.Range("B18").Value = 5437.4721692008497
$ws.Range("C18").Value = 0.27168383974023569
$ws.Range("D18").Value = 114.53607911419446

$ws.Range("B19").Value = 1513.6501832768045
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 115.58767172313233

$ws.Range("B20").Value = 4191.22098253287
$ws.Range("C20").Value = 0.30209940341538127
$ws.Range("D20").Value = 13.90957656407797

$ws.Range("B21").Value = 1212.4220897554305
$ws.Range("C21").Value = 0.15206429208480401
$ws.Range("D21").Value = 186.15377694624723

$ws.Range("B22").Value = 1177.2640827934172
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 102.23883301566036

$ws.Range("B23").Value = 1118.6139578295308
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 129.294262896541

$ws.Range("B24").Value = 1168.6301357915765
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 148.94582454163896

$ws.Range("B25").Value = 1135.6928285138029
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 133.45201597622855

$ws.Range("B26").Value = 3186.3252935344699
$ws.Range("C26").Value = 0.32107660708224106
$ws.Range("D26").Value = 62.455109968727967

$ws.Range("B27").Value = 5277.355990070364
$ws.Range("C27").Value = 0.41716555524737875
$ws.Range("D27").Value = 12.828637181894603

$ws.Range("B28").Value = 5170.9245420112675
$ws.Range("C28").Value = 0.43591084833337845
$ws.Range("D28").Value = 11.82596964423858

# Reset the view: the old scroll position (topLeftCell="A7") and the
# B21:D21 selection from browsing the previous results no longer apply
# to the recomputed sheet, so land back on A1.
$ws.Range("A1").Select() | Out-Null
